$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (copy style from existing header cell H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-11
$data = @(
    @(1, 3),
    @(9, 9),
    @(4, 6),
    @(3, 5),
    @(1, 3),
    @(1, 4),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]   # column I
    $ws.Cells.Item($row, 10).Value = $data[$i][1]  # column J
}
